# The deck currently has zero slides (no <p:sldIdLst> at all). The
# commit "Added 1 File to Branch" adds a single new slide, built from
# the master's Blank layout (slideLayout7.xml - type="blank"), with an
# otherwise empty shape tree - i.e. a brand new, unpopulated slide.
$p = $ppt.ActivePresentation

# ppLayoutBlank = 12
$s = $p.Slides.Add(1, 12)
